$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings
# (e.g. "1.020", "28.384.27") are preserved exactly as text and are
# not reinterpreted as numbers (which would drop trailing zeros /
# convert to floating point).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.384.27"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").Value = "1.869.06"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("D4").Value = "1.020"
$ws.Range("E4").Value = "  +1.82%  "
$ws.Range("D5").Value = "316.57"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("D6").Value = "1.017"
$ws.Range("E6").Value = "  +1.60%  "
$ws.Range("D7").Value = "0.5109"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "0.3955"
$ws.Range("E8").Value = "  +1.94%  "
$ws.Range("D9").Value = "0.08465"
$ws.Range("E9").Value = "  +1.45%  "
$ws.Range("D10").Value = "1.110"
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "41.87"
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "6.245"
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").Value = "1.860.04"
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("B14").Value = "Solana"
$ws.Range("C14").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D14").Value = "20.44"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "7.223"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").Value = "1.019"
$ws.Range("E16").Value = "  +1.73%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.00001112"
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").Value = "90.74"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "0.06771"
$ws.Range("E19").Value = "  +1.77%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "17.67"
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "1.017"
$ws.Range("E21").Value = "  +1.69%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "5.934"
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("B23").Value = "WrappedBTC"
$ws.Range("C23").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D23").Value = "28.382.26"
$ws.Range("E23").Value = "  +1.03%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "11.15"
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "2.280"
$ws.Range("E25").Value = "  +1.33%  "
$ws.Range("D26").Value = "2.068.17"
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "161.81"
$ws.Range("E27").Value = "  +1.83%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "20.77"
$ws.Range("E28").Value = "  +1.13%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "2.361"
$ws.Range("E29").Value = "  -3.73%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "127.26"
$ws.Range("E30").Value = "  +1.98%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.1051"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "1.039"
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "5.772"
$ws.Range("E33").Value = "  -0.98%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "3.645"
$ws.Range("E34").Value = "  +1.41%  "
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").Value = "0.02431"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.06451"
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").Value = "0.2179"
$ws.Range("E37").Value = "  -1.42%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "8.837"
$ws.Range("E38").Value = "  -6.62%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "1.256"
$ws.Range("E39").Value = "  +1.40%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "1.180"
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.6369"
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("B42").Value = "InternetComputer(DFINITY)"
$ws.Range("C42").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D42").Value = "4.986"
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "11.23"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").Value = "0.6008"
$ws.Range("E44").Value = "  -0.97%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "12.91"
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "3.701"
$ws.Range("E46").Value = "  +0.87%  "
$ws.Range("B47").Value = "WEMIXTOKEN"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "1.211"
$ws.Range("E47").Value = "  -4.89%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "1.987"
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").Value = "1.203"
$ws.Range("E49").Value = "  -2.20%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "121.25"
$ws.Range("E50").Value = "  +0.79%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.06854"
$ws.Range("E51").Value = "  -0.40%  "
